$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # "MLX90614 Housing"

# NOTE: the order in which new text values are written below matches the
# order new strings appear in the shared-string table of the target file
# ("1/2" PVC Coupling", then "115991", the Lowe's CARLON link, "23761",
# then the Lowe's Genova link), so set B8 before the rest of row 4/row 8.

# --- Row 8: "1/2" Coupling" -> "1/2" PVC Coupling" ---
$ws.Range("B8").Value = "1/2`" PVC Coupling"

# --- Row 4: new PVC elbow line gets pricing + source info ---
$ws.Range("D4").Value = 2.48
$ws.Range("F4").Formula = "=D4"
$ws.Range("H4").Formula = "=D4"
$ws.Range("J4").Formula = "=D4"
$ws.Range("L4").Value = "Lowe's"
$ws.Range("M4").Value = "115991"
$ws.Range("N4").Value = "N/A"
$ws.Range("O4").Value = "https://www.lowes.com/pd/CARLON-1-2-in-PVC-Transition/3127629"

# --- Row 8 (continued): pricing + new source info ---
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0.33
$ws.Range("F8").Formula = "=D8"
$ws.Range("H8").Formula = "=D8"
$ws.Range("J8").Value = 0.33
$ws.Range("L8").Value = "Lowe's"
$ws.Range("M8").Value = "23761"
$ws.Range("N8").Value = "N/A"
$ws.Range("O8").Value = "https://www.lowes.com/pd/Genova-1-2-in-Dia-Coupling-CPVC-Fittings/1000200923"

# Column L now holds data, matching the author's resulting bestFit width
$ws.Columns.Item(12).ColumnWidth = 10.666666666666666

# --- Window/view state: MLX90614 Housing becomes the active/selected tab ---
$ws.Activate()
$ws.Range("C34").Select()
